$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.859.36"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "1.839.91"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.24"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4713"
$ws.Range("E7").Value = "  +3.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3647"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07144"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.959.40"
$ws.Range("E11").Value = "  +8.22%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.49"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07612"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.276"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.390"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.70"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008621"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "26.905.29"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.46"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.009"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.924"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.55"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.17"
$ws.Range("E26").Value = "  +2.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.005"
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.07"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.852"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08818"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.220"
$ws.Range("E31").Value = "  +3.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.168"
$ws.Range("E32").Value = "  +5.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7420"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.472"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.742"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.088"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01941"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05226"
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.965"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5180"
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.946"
$ws.Range("E41").Value = "  +2.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1509"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.137"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("E44").Value = "  +5.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4694"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.007"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.67"
$ws.Range("E47").Value = "  +2.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.592"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.05"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06028"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8833"
$ws.Range("E51").Value = "  +4.45%  "
